$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The totals that used to sit on row 9 move down to row 14, making room for
# a new attendance entry on row 9 (continuing the B/C/D/E log started above).
$ws.Range("F14").Formula = "=E14*4"
$ws.Range("E14").Formula = "=SUM(E2:E8)"
$ws.Range("E14").NumberFormat = $ws.Range("E9").NumberFormat

# Remove the old row-9 totals entirely (value, formula and formatting).
$ws.Range("E9:F9").Clear()

# New data entry on row 9: a time value (20:00), formatted like the "od"
# column entries above it (e.g. B7).
$ws.Range("B9").Value = 0.83333333333333337
$ws.Range("B9").NumberFormat = $ws.Range("B7").NumberFormat

# Update the active selection to match the new cursor position.
$ws.Range("B10").Select()
